# Update the fixed "Date" placeholder text that appears in the
# Insert > Header & Footer dialog (stored as the <a:fld type="datetime1"/>
# / <a:fld type="datetimeFigureOut"/> placeholder shape on the slide
# master, every slide layout, and the notes master) from 7/20/2019 to
# 9/18/2019.

$p = $ppt.ActivePresentation
$oldDate = "7/20/2019"
$newDate = "9/18/2019"

function Update-DatePlaceholder($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $phType = -1
        try {
            $phType = $shp.PlaceholderFormat.Type
        } catch {
            $phType = -1
        }
        # ppPlaceholderDate = 16
        if ($phType -eq 16) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1. Slide master.
$master = $p.SlideMaster
Update-DatePlaceholder $master

# 2. Every slide layout belonging to the slide master.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li)
}

# 3. Notes master.
Update-DatePlaceholder $p.NotesMaster
